# Add a new weekly Piña price record at row 232 of the "Macroferia Regional
# de Talca" sheet. Inserting a full row pushes the existing rows 232-271
# down to 233-272 (preserving all of their data/formatting), matching the
# target dimension A1:T272. We then populate the newly inserted row 232
# with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 232, shifting rows
# 232:271 down to 233:272.
$ws.Rows(232).Insert()

# Fill in the new record in row 232.
$ws.Range("A232").Value = 5
$ws.Range("B232").Value = 'Macroferia Regional de Talca'
$ws.Range("C232").Value = 'Maule'
$ws.Range("D232").Value = 44776
$ws.Range("E232").Value = 7
$ws.Range("F232").Value = 'Fruta'
$ws.Range("G232").Value = 100108
$ws.Range("H232").Value = 'Tropicales y subtropicales'
$ws.Range("I232").Value = 100108005
$ws.Range("J232").Value = 'Piña'
$ws.Range("K232").Value = 'Caramelo'
$ws.Range("L232").Value = 'Segunda'
$ws.Range("M232").Value = 230
$ws.Range("N232").Value = 19000
$ws.Range("O232").Value = 19000
$ws.Range("P232").Value = 19000
$ws.Range("Q232").Value = '$/caja 14 unidades'
$ws.Range("R232").Value = 'Ecuador'
$ws.Range("S232").Value = 1357
$ws.Range("T232").Value = 14
